# Add a new row of data (an "AttributeBarButtonItem" / Thuoc tinh bien the
# migration entry) to the Icons sheet, as the sheet's table grows by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New data row: MODULE | GUI | PHÂN LOẠI | Control | TEXT / CAPTION | FILE NAME
$ws.Range("A32").Value = "VnsErp2025"
$ws.Range("B32").Value = "FormMain"
$ws.Range("C32").Value = "SanPhamDichVuRibbonPageGroup"
$ws.Range("D32").Value = "AttributeBarButtonItem"
$ws.Range("E32").Value = "Thuộc tính biến thể"
$ws.Range("F32").Value = "data-classification.svg"

# Move the active selection down to the newly-added row, matching the
# author's recorded cursor position after the edit.
$ws.Range("A32:B32").Select()
